# Bugfixed the naive forecaster component module
# The underlying forecast series shifted by one period (oldest row dropped,
# a new row appended) and the y_1_forecast values were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 3-19 (original) up into rows 2-18, carrying over the
# date / y_0 / y_0_forecast / y_1 columns, and write the freshly
# recalculated y_1_forecast (column E) values.

$ws.Range("A2").Value = 39765
$ws.Range("B2").Value = 2008
$ws.Range("D2").Value = 2009
$ws.Range("E2").Value = 3.407109591918855

$ws.Range("A3").Value = 40130
$ws.Range("B3").Value = 2009
$ws.Range("C3").Value = -4.715480642250625
$ws.Range("D3").Value = 2010
$ws.Range("E3").Value = -2.725947775269033

$ws.Range("A4").Value = 40494
$ws.Range("B4").Value = 2010
$ws.Range("C4").Value = 6.130685532900904
$ws.Range("D4").Value = 2011
$ws.Range("E4").Value = 5.993806847197725

$ws.Range("A5").Value = 40862
$ws.Range("B5").Value = 2011
$ws.Range("C5").Value = 8.703939237318981
$ws.Range("D5").Value = 2012
$ws.Range("E5").Value = 6.924353497010971

$ws.Range("A6").Value = 41228
$ws.Range("B6").Value = 2012
$ws.Range("C6").Value = 2.688274587589135
$ws.Range("D6").Value = 2013
$ws.Range("E6").Value = 3.737237443362851

$ws.Range("A7").Value = 41592
$ws.Range("B7").Value = 2013
$ws.Range("C7").Value = 0.9946838291217786
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = 3.504647804006344

$ws.Range("A8").Value = 41957
$ws.Range("B8").Value = 2014
$ws.Range("C8").Value = 5.562499360312567
$ws.Range("D8").Value = 2015
$ws.Range("E8").Value = 4.130094879572455

$ws.Range("A9").Value = 42321
$ws.Range("B9").Value = 2015
$ws.Range("C9").Value = 4.195080504802551
$ws.Range("D9").Value = 2016
$ws.Range("E9").Value = 5.194458387461709

$ws.Range("A10").Value = 42689
$ws.Range("B10").Value = 2016
$ws.Range("C10").Value = 4.230623896992025
$ws.Range("D10").Value = 2017
$ws.Range("E10").Value = 4.372655645302403

$ws.Range("A11").Value = 43053
$ws.Range("B11").Value = 2017
$ws.Range("C11").Value = 4.933871867981643
$ws.Range("D11").Value = 2018
$ws.Range("E11").Value = 5.324897060120137

$ws.Range("A12").Value = 43418
$ws.Range("B12").Value = 2018
$ws.Range("C12").Value = 5.456119081407906
$ws.Range("D12").Value = 2019
$ws.Range("E12").Value = 4.693063499664252

$ws.Range("A13").Value = 43783
$ws.Range("B13").Value = 2019
$ws.Range("C13").Value = 3.346849276607955
$ws.Range("D13").Value = 2020
$ws.Range("E13").Value = 3.273620772016161

$ws.Range("A14").Value = 44159
$ws.Range("B14").Value = 2020
$ws.Range("C14").Value = -9.2489161297999
$ws.Range("D14").Value = 2021
$ws.Range("E14").Value = -4.319815935184923

$ws.Range("A15").Value = 44525
$ws.Range("B15").Value = 2021
$ws.Range("C15").Value = -1.287084480507283
$ws.Range("D15").Value = 2022
$ws.Range("E15").Value = 1.850145327219943

$ws.Range("A16").Value = 44890
$ws.Range("B16").Value = 2022
$ws.Range("C16").Value = 1.494343500592232
$ws.Range("D16").Value = 2023
$ws.Range("E16").Value = 2.552834403233084

$ws.Range("A17").Value = 45254
$ws.Range("B17").Value = 2023
$ws.Range("C17").Value = -0.6982718287330991
$ws.Range("D17").Value = 2024
$ws.Range("E17").Value = -1.04392885455985

$ws.Range("A18").Value = 45618
$ws.Range("B18").Value = 2024
$ws.Range("C18").Value = -0.4137309550271362
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = 0.8860470190541037

# Row 19 (old last row) is no longer part of the series - remove it
# entirely so the sheet's dimension shrinks back to A1:E18.
$ws.Range("A19:E19").ClearContents()
$ws.Rows.Item(19).Delete()
